$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Server TCP")
$cell = $ws.Cells.Item(30, 1)
$v = $cell.Value()
Write-Host "Value: $v"
$v2 = $cell.Value2()
Write-Host "Value2: $v2"
$r = $ws.Range("A30")
Write-Host "Range A30 value: $($r.Value())"
